# swapped bonferroni-correction with fdr-correction:
#   - header cell A3 on both sheets: "targt" -> "target"
#   - specific significance cells flip from "< .05" to "< .05*" (FDR flags
#     additional comparisons as significant compared to the Bonferroni
#     correction) on both the "arousal" and "valence" sheets.

$wb = $excel.ActiveWorkbook

$sheetEdits = @{
    "arousal" = @{
        "A3" = "target"
        "E29" = "< .05*"
        "H29" = "< .05*"
        "E39" = "< .05*"
        "E51" = "< .05*"
        "H51" = "< .05*"
        "E59" = "< .05*"
        "H59" = "< .05*"
        "E61" = "< .05*"
        "H61" = "< .05*"
        "E93" = "< .05*"
        "H93" = "< .05*"
    }
    "valence" = @{
        "A3" = "target"
        "E29" = "< .05*"
        "H29" = "< .05*"
        "E39" = "< .05*"
        "H39" = "< .05*"
        "E82" = "< .05*"
        "E84" = "< .05*"
        "E86" = "< .05*"
        "E93" = "< .05*"
        "H93" = "< .05*"
    }
}

foreach ($sheetName in $sheetEdits.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $sheetEdits[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
